$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217, shifting existing rows 217:258 down to 218:259
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new record
$ws.Cells.Item(217, 1).Value = 11
$ws.Cells.Item(217, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(217, 3).Value = "Bíobío"
$ws.Cells.Item(217, 4).Value = 45015
$ws.Cells.Item(217, 5).Value = 8
$ws.Cells.Item(217, 6).Value = "Fruta"
$ws.Cells.Item(217, 7).Value = 100108
$ws.Cells.Item(217, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(217, 9).Value = 100108005
$ws.Cells.Item(217, 10).Value = "Piña"
$ws.Cells.Item(217, 11).Value = "Caramelo"
$ws.Cells.Item(217, 12).Value = "Primera"
$ws.Cells.Item(217, 13).Value = 100
$ws.Cells.Item(217, 14).Value = 20000
$ws.Cells.Item(217, 15).Value = 21000
$ws.Cells.Item(217, 16).Value = 20500
$ws.Cells.Item(217, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(217, 18).Value = "Ecuador"
$ws.Cells.Item(217, 19).Value = 1708
$ws.Cells.Item(217, 20).Value = 12
